$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474758173/gFRHQATQ91-camera-video-segment-driver-1749474760673.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNA7ZKPHX%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T190004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFqnarxUs0V7zDpi5gegPSsSXkXVR5DxVf1avdpuMzHqAiEAyLxASvhFKLEmYPUIXmD1%2FN2A3GDZHaNsmFn8Y7sjYq8q5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKJQp2bACHci%2FXvWvCq6A1TtKAL9nn6ix%2FwtI4jaxYGEYwNs65BHZJuV8550fXZxQaHSfZJC6lpLtMxN6NUgh0DkrOdyiApwxpIpqo7Mdef9l9agoPrMr4%2FQD3lWxQF1eELlXiY8n%2FgGPMGe%2BXWfh5macjw9VIORbOdbFGziueIEG02rIVtjeJXBnwO%2BtPSkItu1oLzDSa7i%2BOdAuo95fU0z6QF3FHPW3TKs%2BItAS7tRF2zupDsEC3AroDb7Zd9bPJkSY5BfgdT7fQtV3fTJp%2B9Dp%2Bmf7PlClHBmUVktZbS5pOvSo%2F4NQkAv3EEWyrCIQXRAxBLAwzZDG2SVTJ5leDNLhstwE%2BYGNuqZW4wzrXQmx6fHN8XpgFAI0W%2F1JTEN9O30mH3JqPeNoew710qSWLgXWYxDFSeJHLVdMa0q%2FjK%2BOsaq9%2BKx2Eoi0wn4334j0fuXik5fMpOvT%2BYygWcVqhOBbZqdfGuzA%2BerpZrNRLymZcmbBxFPEyU3R5cAKPr5y8iQv%2FMqhmeaPqi8ROsBHut%2BpUOeY5Awva0H4jisR%2BTvo5eAGrJRi1wo3B1GY9PK7wZ5vX9dnAhpY6y1cD6WYDgLFLatFNlx%2BC8wzOihwgY6pQF0%2FtrTrq%2BZQFXRBuPsur51sygeAq%2BHb7hbvexb9%2BaZXIjTArUsxK2hI6jTEZAW5lwuCmFai92Q8nG0rBsCUy0lstjWKbSWns5l74HW6a48zhbfbDcQFS1aAOd0%2FjKzLUzoSax%2FIFNFkxnog78YY0jaagUrjQVXU2%2F1U54V58xrk%2FgL6ocFKmuHYs%2B49RKi0Q%2FhYA98qcx6vA3jFy0Xn0k20Lgsec8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2003%3A00%3A04%20GMT&X-Amz-Signature=4fa5a464e0c19c555caee7f8571f1703438416bace8b9ec4d951d5d3b4933c6a'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474758173/7Z3ova6Wtr-camera-video-segment-1749474760673.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNA7ZKPHX%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T190004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFqnarxUs0V7zDpi5gegPSsSXkXVR5DxVf1avdpuMzHqAiEAyLxASvhFKLEmYPUIXmD1%2FN2A3GDZHaNsmFn8Y7sjYq8q5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKJQp2bACHci%2FXvWvCq6A1TtKAL9nn6ix%2FwtI4jaxYGEYwNs65BHZJuV8550fXZxQaHSfZJC6lpLtMxN6NUgh0DkrOdyiApwxpIpqo7Mdef9l9agoPrMr4%2FQD3lWxQF1eELlXiY8n%2FgGPMGe%2BXWfh5macjw9VIORbOdbFGziueIEG02rIVtjeJXBnwO%2BtPSkItu1oLzDSa7i%2BOdAuo95fU0z6QF3FHPW3TKs%2BItAS7tRF2zupDsEC3AroDb7Zd9bPJkSY5BfgdT7fQtV3fTJp%2B9Dp%2Bmf7PlClHBmUVktZbS5pOvSo%2F4NQkAv3EEWyrCIQXRAxBLAwzZDG2SVTJ5leDNLhstwE%2BYGNuqZW4wzrXQmx6fHN8XpgFAI0W%2F1JTEN9O30mH3JqPeNoew710qSWLgXWYxDFSeJHLVdMa0q%2FjK%2BOsaq9%2BKx2Eoi0wn4334j0fuXik5fMpOvT%2BYygWcVqhOBbZqdfGuzA%2BerpZrNRLymZcmbBxFPEyU3R5cAKPr5y8iQv%2FMqhmeaPqi8ROsBHut%2BpUOeY5Awva0H4jisR%2BTvo5eAGrJRi1wo3B1GY9PK7wZ5vX9dnAhpY6y1cD6WYDgLFLatFNlx%2BC8wzOihwgY6pQF0%2FtrTrq%2BZQFXRBuPsur51sygeAq%2BHb7hbvexb9%2BaZXIjTArUsxK2hI6jTEZAW5lwuCmFai92Q8nG0rBsCUy0lstjWKbSWns5l74HW6a48zhbfbDcQFS1aAOd0%2FjKzLUzoSax%2FIFNFkxnog78YY0jaagUrjQVXU2%2F1U54V58xrk%2FgL6ocFKmuHYs%2B49RKi0Q%2FhYA98qcx6vA3jFy0Xn0k20Lgsec8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2003%3A00%3A04%20GMT&X-Amz-Signature=d08c98a919d59f13012c73546f266c74bfb6caa89b5d50f1c43bc89eb0e9a241'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474688613/JLWgs3viB6-camera-video-segment-driver-1749474691113.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNA7ZKPHX%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T190004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFqnarxUs0V7zDpi5gegPSsSXkXVR5DxVf1avdpuMzHqAiEAyLxASvhFKLEmYPUIXmD1%2FN2A3GDZHaNsmFn8Y7sjYq8q5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKJQp2bACHci%2FXvWvCq6A1TtKAL9nn6ix%2FwtI4jaxYGEYwNs65BHZJuV8550fXZxQaHSfZJC6lpLtMxN6NUgh0DkrOdyiApwxpIpqo7Mdef9l9agoPrMr4%2FQD3lWxQF1eELlXiY8n%2FgGPMGe%2BXWfh5macjw9VIORbOdbFGziueIEG02rIVtjeJXBnwO%2BtPSkItu1oLzDSa7i%2BOdAuo95fU0z6QF3FHPW3TKs%2BItAS7tRF2zupDsEC3AroDb7Zd9bPJkSY5BfgdT7fQtV3fTJp%2B9Dp%2Bmf7PlClHBmUVktZbS5pOvSo%2F4NQkAv3EEWyrCIQXRAxBLAwzZDG2SVTJ5leDNLhstwE%2BYGNuqZW4wzrXQmx6fHN8XpgFAI0W%2F1JTEN9O30mH3JqPeNoew710qSWLgXWYxDFSeJHLVdMa0q%2FjK%2BOsaq9%2BKx2Eoi0wn4334j0fuXik5fMpOvT%2BYygWcVqhOBbZqdfGuzA%2BerpZrNRLymZcmbBxFPEyU3R5cAKPr5y8iQv%2FMqhmeaPqi8ROsBHut%2BpUOeY5Awva0H4jisR%2BTvo5eAGrJRi1wo3B1GY9PK7wZ5vX9dnAhpY6y1cD6WYDgLFLatFNlx%2BC8wzOihwgY6pQF0%2FtrTrq%2BZQFXRBuPsur51sygeAq%2BHb7hbvexb9%2BaZXIjTArUsxK2hI6jTEZAW5lwuCmFai92Q8nG0rBsCUy0lstjWKbSWns5l74HW6a48zhbfbDcQFS1aAOd0%2FjKzLUzoSax%2FIFNFkxnog78YY0jaagUrjQVXU2%2F1U54V58xrk%2FgL6ocFKmuHYs%2B49RKi0Q%2FhYA98qcx6vA3jFy0Xn0k20Lgsec8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2003%3A00%3A04%20GMT&X-Amz-Signature=ba7aa03bbc5c0e3921ab26602f9e7a345c01106721b9b441e877c930416396e3'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474688613/LOHXs2uAgt-camera-video-segment-1749474691113.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSNA7ZKPHX%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T190004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIFqnarxUs0V7zDpi5gegPSsSXkXVR5DxVf1avdpuMzHqAiEAyLxASvhFKLEmYPUIXmD1%2FN2A3GDZHaNsmFn8Y7sjYq8q5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDKJQp2bACHci%2FXvWvCq6A1TtKAL9nn6ix%2FwtI4jaxYGEYwNs65BHZJuV8550fXZxQaHSfZJC6lpLtMxN6NUgh0DkrOdyiApwxpIpqo7Mdef9l9agoPrMr4%2FQD3lWxQF1eELlXiY8n%2FgGPMGe%2BXWfh5macjw9VIORbOdbFGziueIEG02rIVtjeJXBnwO%2BtPSkItu1oLzDSa7i%2BOdAuo95fU0z6QF3FHPW3TKs%2BItAS7tRF2zupDsEC3AroDb7Zd9bPJkSY5BfgdT7fQtV3fTJp%2B9Dp%2Bmf7PlClHBmUVktZbS5pOvSo%2F4NQkAv3EEWyrCIQXRAxBLAwzZDG2SVTJ5leDNLhstwE%2BYGNuqZW4wzrXQmx6fHN8XpgFAI0W%2F1JTEN9O30mH3JqPeNoew710qSWLgXWYxDFSeJHLVdMa0q%2FjK%2BOsaq9%2BKx2Eoi0wn4334j0fuXik5fMpOvT%2BYygWcVqhOBbZqdfGuzA%2BerpZrNRLymZcmbBxFPEyU3R5cAKPr5y8iQv%2FMqhmeaPqi8ROsBHut%2BpUOeY5Awva0H4jisR%2BTvo5eAGrJRi1wo3B1GY9PK7wZ5vX9dnAhpY6y1cD6WYDgLFLatFNlx%2BC8wzOihwgY6pQF0%2FtrTrq%2BZQFXRBuPsur51sygeAq%2BHb7hbvexb9%2BaZXIjTArUsxK2hI6jTEZAW5lwuCmFai92Q8nG0rBsCUy0lstjWKbSWns5l74HW6a48zhbfbDcQFS1aAOd0%2FjKzLUzoSax%2FIFNFkxnog78YY0jaagUrjQVXU2%2F1U54V58xrk%2FgL6ocFKmuHYs%2B49RKi0Q%2FhYA98qcx6vA3jFy0Xn0k20Lgsec8%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2003%3A00%3A04%20GMT&X-Amz-Signature=e921ad7df7f08c20b1c458a9c3b7945d7c412b6701351066101ba699d40d09fd'

Write-Host "done"
